$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '41.760.19'
Set-TextCell 'E2' '  +1.34%  '

Set-TextCell 'D3' '2.270.06'
Set-TextCell 'E3' '  +0.94%  '

Set-TextCell 'E4' '  +0.03%  '

Set-TextCell 'D5' '303.46'
Set-TextCell 'E5' '  +0.23%  '

Set-TextCell 'D6' '92.53'
Set-TextCell 'E6' '  +1.41%  '

Set-TextCell 'D7' '0.529'
Set-TextCell 'E7' '  +1.71%  '

Set-TextCell 'E8' '  -0.04%  '

Set-TextCell 'D9' '0.483'
Set-TextCell 'E9' '  -0.05%  '

Set-TextCell 'D10' '32.50'
Set-TextCell 'E10' '  +1.81%  '

Set-TextCell 'D11' '53.33'
Set-TextCell 'E11' '  -0.67%  '

Set-TextCell 'E12' '  +0.31%  '

Set-TextCell 'E13' '  -1.74%  '

Set-TextCell 'D14' '6.67'
Set-TextCell 'E14' '  +1.37%  '

Set-TextCell 'D15' '2.619.08'
Set-TextCell 'E15' '  +0.80%  '

Set-TextCell 'D16' '14.25'
Set-TextCell 'E16' '  +1.00%  '

Set-TextCell 'D17' '2.262.72'
Set-TextCell 'E17' '  +2.07%  '

Set-TextCell 'D18' '0.774'
Set-TextCell 'E18' '  +3.37%  '

Set-TextCell 'D19' '41.655.95'
Set-TextCell 'E19' '  +1.25%  '

Set-TextCell 'D20' '12.51'
Set-TextCell 'E20' '  +4.19%  '

Set-TextCell 'E21' '  +0.22%  '

Set-TextCell 'D22' '5.94'
Set-TextCell 'E22' '  +1.33%  '

Set-TextCell 'D23' '67.03'
Set-TextCell 'E23' '  +0.53%  '

Set-TextCell 'D24' '239.90'
Set-TextCell 'E24' '  -0.29%  '

Set-TextCell 'E25' '  +0.75%  '

Set-TextCell 'E26' '  -0.02%  '

Set-TextCell 'D27' '1.93'
Set-TextCell 'E27' '  +4.23%  '

Set-TextCell 'D28' '23.99'
Set-TextCell 'E28' '  +1.00%  '

Set-TextCell 'D29' '9.53'
Set-TextCell 'E29' '  -0.82%  '

Set-TextCell 'E30' '  -1.23%  '

Set-TextCell 'D31' '35.46'
Set-TextCell 'E31' '  +6.34%  '

Set-TextCell 'D32' '160.53'
Set-TextCell 'E32' '  +1.35%  '

Set-TextCell 'E33' '  +1.80%  '

Set-TextCell 'E34' '  +0.06%  '

Set-TextCell 'D35' '0.0743'
Set-TextCell 'E35' '  +1.33%  '

Set-TextCell 'E36' '  -0.46%  '

Set-TextCell 'D37' '16.90'
Set-TextCell 'E37' '  +1.29%  '

Set-TextCell 'E38' '  +0.55%  '

Set-TextCell 'E39' '  +2.15%  '

Set-TextCell 'E40' '  +0.94%  '

Set-TextCell 'E41' '  +0.57%  '

Set-TextCell 'D42' '3.91'
Set-TextCell 'E42' '  -0.15%  '

Set-TextCell 'D43' '2.005.05'
Set-TextCell 'E43' '  -3.00%  '

Set-TextCell 'D44' '19.10'
Set-TextCell 'E44' '  -6.00%  '

Set-TextCell 'D45' '0.0281'
Set-TextCell 'E45' '  +1.73%  '

Set-TextCell 'E46' '  +1.32%  '

Set-TextCell 'E47' '  +4.51%  '

Set-TextCell 'D48' '2.88'
Set-TextCell 'E48' '  -2.12%  '

Set-TextCell 'B49' 'MultiversX'
Set-TextCell 'C49' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextCell 'D49' '52.45'
Set-TextCell 'E49' '  +3.29%  '

Set-TextCell 'B50' 'TrustWalletToken'
Set-TextCell 'C50' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D50' '1.15'
Set-TextCell 'E50' '  +0.85%  '

Set-TextCell 'E51' '  +0.97%  '
